$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '91.742.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.117.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '617.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  -3.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.384'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.114.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.739'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.204'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.74'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.542.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.697.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.107.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '445.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('E24').Value = '  -6.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.280.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +25.55%  '
$ws.Range('E31').Value = '  -3.16%  '
$ws.Range('E32').Value = '  -10.78%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.33'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E36').Value = '  -3.34%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.21'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('B38').Value = 'MantraDAO'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '491.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.439'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.698'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('E49').Value = '  -1.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.20%  '
